$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume(1h)) updates per row
# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.629.75"
$ws.Range("E2").Value = "  +3.00%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.724.68"
$ws.Range("E3").Value = "  +3.80%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "526.59"
$ws.Range("E5").Value = "  +1.10%  "

# Row 6 - Solana
$ws.Range("D6").Value = "145.51"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.92%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.725.50"
$ws.Range("E9").Value = "  +3.28%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +7.95%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.16%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.24%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +3.26%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.195.68"
$ws.Range("E14").Value = "  +3.48%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "60.638.84"
$ws.Range("E15").Value = "  +3.06%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "21.27"
$ws.Range("E16").Value = "  +1.93%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.721.85"
$ws.Range("E17").Value = "  +3.32%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +1.00%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "344.55"
$ws.Range("E19").Value = "  -0.15%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +0.49%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "10.60"
$ws.Range("E21").Value = "  +3.78%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.49"
$ws.Range("E22").Value = "  +5.60%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.03%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "63.35"
$ws.Range("E24").Value = "  +2.97%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +1.19%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +2.45%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +0.07%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0" + [char]8323 + "0820"
$ws.Range("E28").Value = "  +2.64%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("E29").Value = "  +2.51%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "6.84"
$ws.Range("E30").Value = "  +9.47%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +1.88%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "18.99"
$ws.Range("E33").Value = "  +0.68%  "

# Row 34 - Monero
$ws.Range("E34").Value = "  -0.47%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +7.26%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +7.76%  "

# Row 37 - SuiNetwork
$ws.Range("D37").Value = "0.942"
$ws.Range("E37").Value = "  -3.66%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "0.875"
$ws.Range("E38").Value = "  +4.18%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +7.18%  "

# Row 40 - OKB
$ws.Range("E40").Value = "  +1.33%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +0.52%  "

# Row 42 - Bittensor
$ws.Range("D42").Value = "280.90"
$ws.Range("E42").Value = "  +1.24%  "

# Row 43 - EnergySwap
$ws.Range("E43").Value = "  +3.31%  "

# Row 44 - FirstDigitalUSD
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.28%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "0.610"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.140.35"
$ws.Range("E46").Value = "  +7.61%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "0.0986"
$ws.Range("E47").Value = "  +0.09%  "

# Row 48 and 49 - swap Hedera and RenderToken entries
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "4.92"
$ws.Range("E48").Value = "  +6.20%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0538"
$ws.Range("E49").Value = "  +2.97%  "

# Row 50 - WhiteBITCoin
$ws.Range("E50").Value = "  +2.23%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  +1.63%  "
